$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    "E2" = 3
    "G2" = 20.32821866666667
    "H2" = 60.984656
    "I2" = 0.004181898474048532
    "J2" = 0.004181898474048532
    "K2" = 3
    "M2" = 153.5290173333333
    "N2" = 460.587052
    "O2" = 0.3172206968818489
    "P2" = 0.317220696881849
    "Q2" = 3120.971436030457
    "R2" = 28088.74292427411
    "S2" = 0.001326584748226816
    "T2" = 0.001326584748226816
    "E3" = 3
    "G3" = 20.32821866666667
    "H3" = 60.984656
    "I3" = 0.004181898474048532
    "J3" = 0.004181898474048532
    "K3" = 3
    "M3" = 168.7997026666667
    "N3" = 506.3991080000001
    "O3" = 0.3487728915577651
    "P3" = 0.3487728915577651
    "Q3" = 3431.397266676317
    "R3" = 30882.57540008685
    "S3" = 0.001458532822994912
    "T3" = 0.001458532822994912
    "E4" = 3
    "G4" = 20.32821866666667
    "H4" = 60.984656
    "I4" = 0.004181898474048532
    "J4" = 0.004181898474048532
    "K4" = 3
    "M4" = 68.09032333333333
    "N4" = 204.27097
    "O4" = 0.1406878008722904
    "P4" = 0.1406878008722904
    "Q4" = 1384.154981804036
    "R4" = 12457.39483623632
    "S4" = 0.0005883420997850749
    "T4" = 0.000588342099785075
    "E5" = 3
    "G5" = 20.32821866666667
    "H5" = 60.984656
    "I5" = 0.004181898474048532
    "J5" = 0.004181898474048532
    "K5" = 3
    "M5" = 93.562673
    "N5" = 280.688019
    "O5" = 0.1933186106880956
    "P5" = 0.1933186106880956
    "Q5" = 1901.96247578183
    "R5" = 17117.66228203646
    "S5" = 0.0008084388030417294
    "T5" = 0.0008084388030417294
    "E6" = 3
    "G6" = 4809.896321333334
    "H6" = 14429.688964
    "I6" = 0.9894865072215304
    "J6" = 0.9894865072215304
    "K6" = 3
    "M6" = 153.5290173333333
    "N6" = 460.587052
    "O6" = 0.3172206968818489
    "P6" = 0.317220696881849
    "Q6" = 738458.6556895216
    "R6" = 6646127.901205694
    "S6" = 0.3138855993760005
    "T6" = 0.3138855993760005
    "E7" = 3
    "G7" = 4809.896321333334
    "H7" = 14429.688964
    "I7" = 0.9894865072215304
    "J7" = 0.9894865072215304
    "K7" = 3
    "M7" = 168.7997026666667
    "N7" = 506.3991080000001
    "O7" = 0.3487728915577651
    "P7" = 0.3487728915577651
    "Q7" = 811909.0688985606
    "R7" = 7307181.620087045
    "S7" = 0.3451060702810466
    "T7" = 0.3451060702810466
    "E8" = 3
    "G8" = 4809.896321333334
    "H8" = 14429.688964
    "I8" = 0.9894865072215304
    "J8" = 0.9894865072215304
    "K8" = 3
    "M8" = 68.09032333333333
    "N8" = 204.27097
    "O8" = 0.1406878008722904
    "P8" = 0.1406878008722904
    "Q8" = 327507.3957193972
    "R8" = 2947566.561474575
    "S8" = 0.1392086806938008
    "T8" = 0.1392086806938008
    "E9" = 3
    "G9" = 4809.896321333334
    "H9" = 14429.688964
    "I9" = 0.9894865072215304
    "J9" = 0.9894865072215304
    "K9" = 3
    "M9" = 93.562673
    "N9" = 280.688019
    "O9" = 0.1933186106880956
    "P9" = 0.1933186106880956
    "Q9" = 450026.7566768136
    "R9" = 4050240.810091322
    "S9" = 0.1912861568706826
    "T9" = 0.1912861568706826
    "E10" = 3
    "G10" = 2.69506
    "H10" = 8.085180000000001
    "I10" = 0.000554424737665286
    "J10" = 0.000554424737665286
    "K10" = 3
    "M10" = 153.5290173333333
    "N10" = 460.587052
    "O10" = 0.3172206968818489
    "P10" = 0.317220696881849
    "Q10" = 413.7699134543733
    "R10" = 3723.92922108936
    "S10" = 0.0001758750016507183
    "T10" = 0.0001758750016507183
    "E11" = 3
    "G11" = 2.69506
    "H11" = 8.085180000000001
    "I11" = 0.000554424737665286
    "J11" = 0.000554424737665286
    "K11" = 3
    "M11" = 168.7997026666667
    "N11" = 506.3991080000001
    "O11" = 0.3487728915577651
    "P11" = 0.3487728915577651
    "Q11" = 454.9253266688268
    "R11" = 4094.327940019441
    "S11" = 0.0001933683189066772
    "T11" = 0.0001933683189066772
    "E12" = 3
    "G12" = 2.69506
    "H12" = 8.085180000000001
    "I12" = 0.000554424737665286
    "J12" = 0.000554424737665286
    "K12" = 3
    "M12" = 68.09032333333333
    "N12" = 204.27097
    "O12" = 0.1406878008722904
    "P12" = 0.1406878008722904
    "Q12" = 183.5075068027333
    "R12" = 1651.5675612246
    "S12" = 0.0000780007970913256
    "T12" = 0.00007800079709132561
    "E13" = 3
    "G13" = 2.69506
    "H13" = 8.085180000000001
    "I13" = 0.000554424737665286
    "J13" = 0.000554424737665286
    "K13" = 3
    "M13" = 93.562673
    "N13" = 280.688019
    "O13" = 0.1933186106880956
    "P13" = 0.1933186106880956
    "Q13" = 252.15701749538
    "R13" = 2269.41315745842
    "S13" = 0.000107180620016565
    "T13" = 0.000107180620016565
    "E14" = 3
    "G14" = 28.08283533333333
    "H14" = 84.24850599999999
    "I14" = 0.005777169566755752
    "J14" = 0.005777169566755752
    "K14" = 3
    "M14" = 153.5290173333333
    "N14" = 460.587052
    "O14" = 0.3172206968818489
    "P14" = 0.317220696881849
    "Q14" = 4311.530112660478
    "R14" = 38803.77101394431
    "S14" = 0.001832637755970869
    "T14" = 0.001832637755970869
    "E15" = 3
    "G15" = 28.08283533333333
    "H15" = 84.24850599999999
    "I15" = 0.005777169566755752
    "J15" = 0.005777169566755752
    "K15" = 3
    "M15" = 168.7997026666667
    "N15" = 506.3991080000001
    "O15" = 0.3487728915577651
    "P15" = 0.3487728915577651
    "Q15" = 4740.374254303628
    "R15" = 42663.36828873265
    "S15" = 0.002014920134816924
    "T15" = 0.002014920134816924
    "E16" = 3
    "G16" = 28.08283533333333
    "H16" = 84.24850599999999
    "I16" = 0.005777169566755752
    "J16" = 0.005777169566755752
    "K16" = 3
    "M16" = 68.09032333333333
    "N16" = 204.27097
    "O16" = 0.1406878008722904
    "P16" = 0.1406878008722904
    "Q16" = 1912.169337963424
    "R16" = 17209.52404167082
    "S16" = 0.0008127772816131894
    "T16" = 0.0008127772816131895
    "E17" = 3
    "G17" = 28.08283533333333
    "H17" = 84.24850599999999
    "I17" = 0.005777169566755752
    "J17" = 0.005777169566755752
    "K17" = 3
    "M17" = 93.562673
    "N17" = 280.688019
    "O17" = 0.1933186106880956
    "P17" = 0.1933186106880956
    "Q17" = 2627.505139205512
    "R17" = 23647.54625284961
    "S17" = 0.001116834394354769
    "T17" = 0.001116834394354769
}

foreach ($key in $changes.Keys) {
    $ws.Range($key).Value = $changes[$key]
}